$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price data scraped as text (e.g. "1.835.84" using dots as
# thousands separators). Pre-format the whole data range as Text so that
# numeric-looking values (e.g. "1.005") are not auto-converted to numbers,
# preserving the original inline-string semantics.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range('D2').Value = '27.469.97'
$ws.Range('E2').Value = '  -1.08%  '
$ws.Range('D3').Value = '1.834.25'
$ws.Range('E3').Value = '  -1.02%  '
$ws.Range('D4').Value = '1.005'
$ws.Range('E4').Value = '  -2.77%  '
$ws.Range('D5').Value = '315.38'
$ws.Range('E5').Value = '  -2.13%  '
$ws.Range('E6').Value = '  -2.38%  '
$ws.Range('D7').Value = '0.4298'
$ws.Range('E7').Value = '  -2.17%  '
$ws.Range('E8').Value = '  -2.53%  '
$ws.Range('D9').Value = '0.07275'
$ws.Range('E9').Value = '  -1.90%  '
$ws.Range('D10').Value = '0.8671'
$ws.Range('E10').Value = '  -2.00%  '
$ws.Range('D11').Value = '21.20'
$ws.Range('E11').Value = '  -1.67%  '
$ws.Range('D12').Value = '1.844.55'
$ws.Range('E12').Value = '  -0.42%  '
$ws.Range('D13').Value = '6.700'
$ws.Range('E13').Value = '  +0.08%  '
$ws.Range('D14').Value = '5.362'
$ws.Range('E14').Value = '  -2.74%  '
$ws.Range('D15').Value = '0.07076'
$ws.Range('E15').Value = '  -1.32%  '
$ws.Range('E16').Value = '  +3.67%  '
$ws.Range('E17').Value = '  -2.73%  '
$ws.Range('D18').Value = '0.000008929'
$ws.Range('E18').Value = '  -1.65%  '
$ws.Range('E19').Value = '  -2.46%  '
$ws.Range('D20').Value = '15.28'
$ws.Range('E20').Value = '  -1.39%  '
$ws.Range('D21').Value = '27.479.04'
$ws.Range('E21').Value = '  -1.05%  '
$ws.Range('D22').Value = '5.172'
$ws.Range('E22').Value = '  -2.01%  '
$ws.Range('D23').Value = '10.95'
$ws.Range('E23').Value = '  -2.90%  '
$ws.Range('D24').Value = '2.064.23'
$ws.Range('E24').Value = '  -1.16%  '
$ws.Range('D25').Value = '2.002'
$ws.Range('E25').Value = '  -2.16%  '
$ws.Range('D26').Value = '153.83'
$ws.Range('E26').Value = '  -3.08%  '
$ws.Range('D27').Value = '18.46'
$ws.Range('E27').Value = '  -1.28%  '
$ws.Range('D28').Value = '2.151'
$ws.Range('E28').Value = '  +7.89%  '
$ws.Range('D29').Value = '5.285'
$ws.Range('E29').Value = '  -1.08%  '
$ws.Range('D30').Value = '117.31'
$ws.Range('E30').Value = '  -0.54%  '
$ws.Range('D31').Value = '0.08865'
$ws.Range('E31').Value = '  -2.22%  '
$ws.Range('D32').Value = '1.210'
$ws.Range('E32').Value = '  -0.40%  '
$ws.Range('D33').Value = '0.7685'
$ws.Range('E33').Value = '  -0.77%  '
$ws.Range('D34').Value = '4.485'
$ws.Range('E34').Value = '  -1.96%  '
$ws.Range('D35').Value = '2.913'
$ws.Range('E35').Value = '  -3.17%  '
$ws.Range('D36').Value = '1.005'
$ws.Range('E36').Value = '  -2.56%  '
$ws.Range('D37').Value = '1.124'
$ws.Range('E37').Value = '  -2.22%  '
$ws.Range('D38').Value = '0.01962'
$ws.Range('D39').Value = '0.05295'
$ws.Range('E39').Value = '  +0.29%  '
$ws.Range('E40').Value = '  +4.37%  '
$ws.Range('D41').Value = '2.876'
$ws.Range('E41').Value = '  +0.99%  '
$ws.Range('E42').Value = '  +0.09%  '
$ws.Range('D43').Value = '0.5087'
$ws.Range('E43').Value = '  -1.89%  '
$ws.Range('D44').Value = '8.690'
$ws.Range('E44').Value = '  -0.56%  '
$ws.Range('D45').Value = '10.60'
$ws.Range('E45').Value = '  -1.17%  '
$ws.Range('B46').Value = 'Decentraland'
$ws.Range('C46').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('D46').Value = '0.4748'
$ws.Range('E46').Value = '  +0.95%  '
$ws.Range('B47').Value = 'Quant'
$ws.Range('C47').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D47').Value = '106.43'
$ws.Range('E47').Value = '  -3.56%  '
$ws.Range('D48').Value = '0.06431'
$ws.Range('E48').Value = '  -2.28%  '
$ws.Range('E49').Value = '  -2.71%  '
$ws.Range('E50').Value = '  -2.28%  '
$ws.Range('D51').Value = '1.833'
$ws.Range('E51').Value = '  -3.13%  '
